$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ManageNewsPage"

# Populate the data for the new sheet.
$newSheet.Range("A1").Value = "Sample test news data"

# Match column width / selection as per the authored sheet.
$newSheet.Columns("A").ColumnWidth = 20.66
$newSheet.Range("B5").Select() | Out-Null

# Make the new sheet the active tab.
$newSheet.Activate() | Out-Null
